# Append the latest stock-ticker refresh to the bottom of the list.
# Data update as of 2023-06-30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tickers = @(
    "AAF", "ABDN", "ABF", "ANTO", "AUTO", "AV", "BARC", "BATS", "BDEV", "BEZ",
    "BF.B", "BKG", "BNZL", "BRBY", "BRK.B", "BT-A", "CCH", "CRDA", "DCC", "DGE",
    "ENT", "EXPN", "FCIT", "FRAS", "GLEN", "HLMA", "HSBA", "HSX", "IMB", "IMI",
    "INF", "ITRK", "JMAT", "KGF", "LGEN", "LLOY", "LSEG", "MNDI", "MNG", "OCDO",
    "PHNX", "PSON", "REL", "RMV", "RR", "RS1", "SBRY", "SDR", "SGRO", "SKG",
    "SMDS", "SMT", "SN", "SPX", "SSE", "STAN", "STJ", "SVT", "ULVR", "UU",
    "WEIR", "WTB"
)

# Find the first empty row at the bottom of column A (currently row 5856,
# right after the existing data ending at row 5855).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($ws.Cells.Item(1, 1).Value -eq $null) {
    $startRow = 1
} else {
    $startRow = $lastRow + 1
}

for ($i = 0; $i -lt $tickers.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $tickers[$i]
}
